$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") values for rows 2-14 were all 45243 (2023-11-13)
# and should become 45244 (2023-11-14).
for ($row = 2; $row -le 14; $row++) {
    $ws.Cells.Item($row, 3).Value = 45244
}
